$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 3 de Octubre de 2020 a las 21:06"

# --- Row 25 / 26: Indonesia and Alemania swap ranking, Alemania gets refreshed stats ---
$ws.Range("A25").Value = "Alemania"
$ws.Range("B25").Value = 299533
$ws.Range("C25").Value = 1170
$ws.Range("D25").Value = 259500
$ws.Range("E25").Value = 30436
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 9597

$ws.Range("A26").Value = "Indonesia"
$ws.Range("B26").Value = 299506
$ws.Range("C26").Value = 4007
$ws.Range("D26").Value = 225052
$ws.Range("E26").Value = 63399
$ws.Range("G26").Value = 83
$ws.Range("H26").Value = 11055

# --- Row 68 / 69: Libano and Paraguay swap ranking, Libano gets refreshed stats ---
$ws.Range("A68").Value = "Libano"
$ws.Range("B68").Value = 43494
$ws.Range("C68").Value = 1321
$ws.Range("D68").Value = 19259
$ws.Range("E68").Value = 23837
$ws.Range("G68").Value = 12
$ws.Range("H68").Value = 398

$ws.Range("A69").Value = "Paraguay"
$ws.Range("B69").Value = 42684
$ws.Range("C69").Value = 0
$ws.Range("D69").Value = 25803
$ws.Range("E69").Value = 15991
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 890

# --- Row 183 / 184: Mauricio and Eritrea swap ranking, Eritrea gets refreshed stats ---
$ws.Range("A183").Value = "Eritrea"
$ws.Range("B183").Value = 398
$ws.Range("C183").Value = 17
$ws.Range("D183").Value = 358
$ws.Range("E183").Value = 40
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 0

$ws.Range("A184").Value = "Mauricio"
$ws.Range("B184").Value = 385
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 344
$ws.Range("E184").Value = 31
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 10

# --- Row 207 / 208: Santa Lucia and Nueva Caledonia swap ranking (tie, no numeric change) ---
$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("B207").Value = 27
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 27
$ws.Range("E207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 0

$ws.Range("A208").Value = "Santa Lucia"
$ws.Range("B208").Value = 27
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 27
$ws.Range("E208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0

# --- Simple in-place numeric refreshes (no country / ranking change) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 7574651
$ws.Range("C4").Value = 25328
$ws.Range("D4").Value = 4798407
$ws.Range("E4").Value = 2562383
$ws.Range("G4").Value = 337
$ws.Range("H4").Value = 213861

# Row 6: Brasil
$ws.Range("B6").Value = 4888926
$ws.Range("C6").Value = 6695
$ws.Range("E6").Value = 510778
$ws.Range("G6").Value = 124
$ws.Range("H6").Value = 145555

# Row 14: Francia
$ws.Range("B14").Value = 606625
$ws.Range("C14").Value = 16972
$ws.Range("E14").Value = 476692

# Row 34: Marruecos
$ws.Range("B34").Value = 131228
$ws.Range("C34").Value = 2663
$ws.Range("D34").Value = 108687
$ws.Range("E34").Value = 20248
$ws.Range("G34").Value = 30
$ws.Range("H34").Value = 2293

# Row 101: Namibia
$ws.Range("B101").Value = 11572
$ws.Range("C101").Value = 92
$ws.Range("D101").Value = 9419
$ws.Range("E101").Value = 2030

# Row 104: Maldivas
$ws.Range("B104").Value = 10465
$ws.Range("C104").Value = 67
$ws.Range("D104").Value = 9310
$ws.Range("E104").Value = 1121

# Row 106: Guayana Francesa
$ws.Range("B106").Value = 10029
$ws.Range("C106").Value = 61
$ws.Range("D106").Value = 9665
$ws.Range("E106").Value = 297

# Row 108: Mozambique
$ws.Range("B108").Value = 9049
$ws.Range("C108").Value = 70
$ws.Range("D108").Value = 5736
$ws.Range("E108").Value = 3249

# Row 118: Cabo Verde
$ws.Range("B118").Value = 6296
$ws.Range("C118").Value = 91
$ws.Range("D118").Value = 5409
$ws.Range("E118").Value = 825

# Row 137: Aruba
$ws.Range("B137").Value = 4074
$ws.Range("C137").Value = 36
$ws.Range("D137").Value = 3474
$ws.Range("E137").Value = 570
